$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fix typo: "По описанием" -> "Под описанием"
$ws.Range("B14").Value = "Под описанием  библиотеки ""Capybara"" Должно быть 3 кнопки "

# Fix typo: "симмволе" -> "символе"
$ws.Range("B68").Value = "Гиперссылка на символе ""Elabs"" в футере ведет на сайт ""http://www.elabs.se/"""

# Fix typo: "гсимвол" -> "символ"
$ws.Range("B70").Value = "1) В копирайте в футере найти символ Elabs`n2) Нажать на символ"

# Update the active selection on the sheet view
$ws.Range("E69").Select()
